$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet data lives inside an Excel Table ("Table3"). Adding a
# ListRow extends the table's range (and autofilter) automatically,
# matching the A1:E50 -> A1:E51 growth in the diff.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# 1 Mayis 2020 (serial 43952) data.
$ws.Range("A51").Value = 43952
$ws.Range("B51").Value = 41.430999999999997
$ws.Range("C51").Value = 2.1880000000000002
$ws.Range("D51").Value = 84
$ws.Range("E51").Value = 4922

# Move the view roughly where the author left it (selection on D57,
# scrolled so row 23 is at the top).
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D57").Select()
